# Auto-update draw results: append the newest Pick 4 draw as a new row
# at the bottom of the "Results" table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Find the first empty row right after the existing data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$newRow = $lastRow + 1

$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
# Keep every field as literal text (matches the rest of the table, which
# stores dates/phase codes/results as text rather than numbers/dates).
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-06"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251006"
$ws.Cells.Item($newRow, 4).Value = "5-5-6-3"
$ws.Cells.Item($newRow, 5).Value = "2025-10-06T21:37:29.870+04:00"
